$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal-looking numeric strings (e.g. "7.50", "1.00")
# must be protected from Excel's automatic text-to-number conversion so they are
# written back as text (matching the workbook's inlineStr cells), preserving exact
# formatting such as trailing zeros.

$ws.Range('D2').Value = '60.764.14'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '3.382.25'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.381.64'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '3.962.21'
$ws.Range('E13').Value = '  -2.28%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.124'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.07'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.401.63'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000171'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.81%  '
$ws.Range('D18').Value = '60.908.36'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('E21').Value = '  -5.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.52'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -3.66%  '
$ws.Range('D27').Value = '3.524.16'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.15%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.09%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.76'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('E36').Value = '  -1.93%  '
$ws.Range('D37').Value = '3.412.85'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '167.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  -3.91%  '
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0778'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.39%  '
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('D48').Value = '2.544.01'
$ws.Range('E48').Value = '  -2.28%  '
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.49%  '
